$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated AE values (column 31) for rows 2-6
$ws.Range("AE2").Value = 52004.8
$ws.Range("AE3").Value = 7050
$ws.Range("AE4").Value = 3377
$ws.Range("AE5").Value = 2015
$ws.Range("AE6").Value = 64446.8

# Updated AG (total) values (column 33) for rows 2-6, reflecting the new row totals
$ws.Range("AG2").Value = 295945.38
$ws.Range("AG3").Value = 131571.81
$ws.Range("AG4").Value = 88013.89999999999
$ws.Range("AG5").Value = 71854.28999999999
$ws.Range("AG6").Value = 587385.38
